$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 21.267222
$ws.Range("H2").Value = 63.801666
$ws.Range("I2").Value = 0.06271644651145813
$ws.Range("J2").Value = 0.06271644651145813
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 3098.650082391316
$ws.Range("R2").Value = 27887.85074152185
$ws.Range("S2").Value = 0.01797414316096817
$ws.Range("T2").Value = 0.01797414316096817
$ws.Range("G3").Value = 21.267222
$ws.Range("H3").Value = 63.801666
$ws.Range("I3").Value = 0.06271644651145813
$ws.Range("J3").Value = 0.06271644651145813
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 3589.900750145992
$ws.Range("R3").Value = 32309.10675131393
$ws.Range("S3").Value = 0.02082370977719273
$ws.Range("T3").Value = 0.02082370977719273
$ws.Range("G4").Value = 21.267222
$ws.Range("H4").Value = 63.801666
$ws.Range("I4").Value = 0.06271644651145813
$ws.Range("J4").Value = 0.06271644651145813
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 2724.887375302336
$ws.Range("R4").Value = 24523.98637772102
$ws.Range("S4").Value = 0.01580608151256681
$ws.Range("T4").Value = 0.01580608151256681
$ws.Range("G5").Value = 21.267222
$ws.Range("H5").Value = 63.801666
$ws.Range("I5").Value = 0.06271644651145813
$ws.Range("J5").Value = 0.06271644651145813
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 1398.555466052538
$ws.Range("R5").Value = 12586.99919447284
$ws.Range("S5").Value = 0.008112512060730426
$ws.Range("T5").Value = 0.008112512060730426
$ws.Range("I6").Value = 0.4054090708715844
$ws.Range("J6").Value = 0.4054090708715843
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 20030.16625996043
$ws.Range("R6").Value = 180271.4963396439
$ws.Range("S6").Value = 0.1161877160446209
$ws.Range("T6").Value = 0.1161877160446208
$ws.Range("I7").Value = 0.4054090708715844
$ws.Range("J7").Value = 0.4054090708715843
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.1346077672198612
$ws.Range("T7").Value = 0.1346077672198612
$ws.Range("I8").Value = 0.4054090708715844
$ws.Range("J8").Value = 0.4054090708715843
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 17614.10476037104
$ws.Range("R8").Value = 158526.9428433393
$ws.Range("S8").Value = 0.1021730212179596
$ws.Range("T8").Value = 0.1021730212179595
$ws.Range("I9").Value = 0.4054090708715844
$ws.Range("J9").Value = 0.4054090708715843
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 9040.484651042019
$ws.Range("R9").Value = 81364.36185937817
$ws.Range("S9").Value = 0.05244056638914282
$ws.Range("T9").Value = 0.05244056638914281
$ws.Range("G10").Value = 121.820091
$ws.Range("H10").Value = 365.460273
$ws.Range("I10").Value = 0.3592440621169263
$ws.Range("J10").Value = 0.3592440621169263
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 17749.2779740611
$ws.Range("R10").Value = 159743.5017665499
$ws.Range("S10").Value = 0.1029571119122267
$ws.Range("T10").Value = 0.1029571119122267
$ws.Range("G11").Value = 121.820091
$ws.Range("H11").Value = 365.460273
$ws.Range("I11").Value = 0.3592440621169263
$ws.Range("J11").Value = 0.3592440621169263
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 20563.19513962628
$ws.Range("R11").Value = 185068.7562566365
$ws.Range("S11").Value = 0.1192796228870516
$ws.Range("T11").Value = 0.1192796228870516
$ws.Range("G12").Value = 121.820091
$ws.Range("H12").Value = 365.460273
$ws.Range("I12").Value = 0.3592440621169263
$ws.Range("J12").Value = 0.3592440621169263
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 15608.33982097341
$ws.Range("R12").Value = 140475.0583887607
$ws.Range("S12").Value = 0.09053830764611882
$ws.Range("T12").Value = 0.09053830764611882
$ws.Range("G13").Value = 121.820091
$ws.Range("H13").Value = 365.460273
$ws.Range("I13").Value = 0.3592440621169263
$ws.Range("J13").Value = 0.3592440621169263
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 8011.020628038189
$ws.Range("R13").Value = 72099.1856523437
$ws.Range("S13").Value = 0.04646901967152917
$ws.Range("T13").Value = 0.04646901967152917
$ws.Range("G14").Value = 58.539182
$ws.Range("H14").Value = 175.617546
$ws.Range("I14").Value = 0.1726304205000311
$ws.Range("J14").Value = 0.1726304205000311
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 8529.20240957753
$ws.Range("R14").Value = 76762.82168619777
$ws.Range("S14").Value = 0.04947480389276847
$ws.Range("T14").Value = 0.04947480389276847
$ws.Range("G15").Value = 58.539182
$ws.Range("H15").Value = 175.617546
$ws.Range("I15").Value = 0.1726304205000311
$ws.Range("J15").Value = 0.1726304205000311
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 9881.396515949888
$ws.Range("R15").Value = 88932.56864354898
$ws.Range("S15").Value = 0.05731839055247859
$ws.Range("T15").Value = 0.05731839055247857
$ws.Range("G16").Value = 58.539182
$ws.Range("H16").Value = 175.617546
$ws.Range("I16").Value = 0.1726304205000311
$ws.Range("J16").Value = 0.1726304205000311
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 7500.400286992151
$ws.Range("R16").Value = 67503.60258292935
$ws.Range("S16").Value = 0.04350709661896527
$ws.Range("T16").Value = 0.04350709661896526
$ws.Range("G17").Value = 58.539182
$ws.Range("H17").Value = 175.617546
$ws.Range("I17").Value = 0.1726304205000311
$ws.Range("J17").Value = 0.1726304205000311
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 3849.599772097378
$ws.Range("R17").Value = 34646.3979488764
$ws.Range("S17").Value = 0.02233012943581882
$ws.Range("T17").Value = 0.02233012943581881
